# Updates currentAveragePrice / LevePrice / LeveProfit data cells (columns H:N)
# across the Leve profit sheets, per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 83.5
$ws.Range("I2").Value = 83.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 83.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 29.5
$ws.Range("N2").ClearContents()
# Row 33
$ws.Range("H33").Value = 621.8
$ws.Range("I33").Value = 543.2941
$ws.Range("K33").Value = 543.2941
$ws.Range("M33").Value = -314.2941
# Row 38
$ws.Range("H38").Value = 549.46155
$ws.Range("I38").Value = 240.5
$ws.Range("J38").Value = 814.2857
$ws.Range("K38").Value = 721.5
$ws.Range("L38").Value = 2442.8571
$ws.Range("M38").Value = -349.5
$ws.Range("N38").Value = -3186.8571
# Row 40
$ws.Range("H40").Value = 1082.0714
$ws.Range("I40").Value = 1001.7778
$ws.Range("J40").Value = 1226.6
$ws.Range("K40").Value = 1001.7778
$ws.Range("L40").Value = 1226.6
$ws.Range("M40").Value = -826.7778
$ws.Range("N40").Value = -1576.6
# Row 43
$ws.Range("H43").Value = 1623.8182
$ws.Range("J43").Value = 1596.8889
$ws.Range("L43").Value = 1596.8889
$ws.Range("N43").Value = -1734.8889
# Row 58
$ws.Range("H58").Value = 2324
$ws.Range("I58").Value = 223.33333
$ws.Range("J58").Value = 5475
$ws.Range("K58").Value = 669.99999
$ws.Range("L58").Value = 16425
$ws.Range("M58").Value = -519.99999
$ws.Range("N58").Value = -16725
# Row 62
$ws.Range("H62").Value = 2737
$ws.Range("I62").Value = 2968.3333
$ws.Range("J62").Value = 2390
$ws.Range("K62").Value = 2968.3333
$ws.Range("L62").Value = 2390
$ws.Range("M62").Value = -2344.3333
$ws.Range("N62").Value = -3638
# Row 64
$ws.Range("H64").Value = 4152.25
$ws.Range("I64").Value = 4043.7
$ws.Range("J64").Value = 4333.1665
$ws.Range("K64").Value = 4043.7
$ws.Range("L64").Value = 4333.1665
$ws.Range("M64").Value = -3795.7
$ws.Range("N64").Value = -4829.1665
# Row 65
$ws.Range("H65").Value = 2737
$ws.Range("I65").Value = 2968.3333
$ws.Range("J65").Value = 2390
$ws.Range("K65").Value = 14841.6665
$ws.Range("L65").Value = 11950
$ws.Range("M65").Value = -11721.6665
$ws.Range("N65").Value = -18190
# Row 67
$ws.Range("H67").Value = 4152.25
$ws.Range("I67").Value = 4043.7
$ws.Range("J67").Value = 4333.1665
$ws.Range("K67").Value = 4043.7
$ws.Range("L67").Value = 4333.1665
$ws.Range("M67").Value = -3185.7
$ws.Range("N67").Value = -6049.1665
# Row 87
$ws.Range("H87").Value = 38639.715
$ws.Range("J87").Value = 38639.715
$ws.Range("L87").Value = 38639.715
$ws.Range("N87").Value = -41135.715
# Row 90
$ws.Range("H90").Value = 38639.715
$ws.Range("J90").Value = 38639.715
$ws.Range("L90").Value = 115919.145
$ws.Range("N90").Value = -128399.145
# Row 113
$ws.Range("H113").Value = 76926490
$ws.Range("I113").Value = 90911490
$ws.Range("K113").Value = 90911490
$ws.Range("M113").Value = -90908236
# Row 137
$ws.Range("H137").Value = 78994.234
$ws.Range("I137").Value = 1880.4
$ws.Range("K137").Value = 5641.200000000001
$ws.Range("M137").Value = -3091.200000000001
# Row 138
$ws.Range("H138").Value = 1840.0927
$ws.Range("I138").Value = 553.61536
$ws.Range("J138").Value = 3034.6785
$ws.Range("K138").Value = 1660.84608
$ws.Range("L138").Value = 9104.0355
$ws.Range("M138").Value = 3479.15392
$ws.Range("N138").Value = -19384.0355

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 42124.27
$ws.Range("I32").Value = 46536.13
$ws.Range("K32").Value = 46536.13
$ws.Range("M32").Value = -46249.13
# Row 45
$ws.Range("H45").Value = 3300.1555
$ws.Range("I45").Value = 3016.3333
$ws.Range("J45").Value = 3489.3704
$ws.Range("K45").Value = 3016.3333
$ws.Range("L45").Value = 3489.3704
$ws.Range("M45").Value = -2639.3333
$ws.Range("N45").Value = -4243.3704
# Row 61
$ws.Range("H61").Value = 2172.16
$ws.Range("I61").Value = 1305.4736
$ws.Range("J61").Value = 4916.6665
$ws.Range("K61").Value = 1305.4736
$ws.Range("L61").Value = 4916.6665
$ws.Range("M61").Value = -1093.4736
$ws.Range("N61").Value = -5340.6665
# Row 132
$ws.Range("H132").Value = 14030.22
$ws.Range("I132").Value = 1452.7576
$ws.Range("J132").Value = 65912.25
$ws.Range("K132").Value = 4358.2728
$ws.Range("L132").Value = 197736.75
$ws.Range("M132").Value = -1828.2728
$ws.Range("N132").Value = -202796.75
# Row 136
$ws.Range("H136").Value = 2172.16
$ws.Range("I136").Value = 1305.4736
$ws.Range("J136").Value = 4916.6665
$ws.Range("K136").Value = 3916.4208
$ws.Range("L136").Value = 14749.9995
$ws.Range("M136").Value = -1366.4208
$ws.Range("N136").Value = -19849.9995

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 248.11111
$ws.Range("I22").Value = 248.11111
$ws.Range("K22").Value = 248.11111
$ws.Range("M22").Value = -75.11111

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2496.0645
$ws.Range("I31").Value = 2105.0952
$ws.Range("K31").Value = 2105.0952
$ws.Range("M31").Value = -1810.0952
# Row 34
$ws.Range("H34").Value = 2496.0645
$ws.Range("I34").Value = 2105.0952
$ws.Range("K34").Value = 2105.0952
$ws.Range("M34").Value = -1903.0952
# Row 58
$ws.Range("H58").Value = 22215.416
$ws.Range("I58").Value = 1870.3
$ws.Range("J58").Value = 36747.645
$ws.Range("K58").Value = 1870.3
$ws.Range("L58").Value = 36747.645
$ws.Range("M58").Value = -1667.3
$ws.Range("N58").Value = -37153.645
# Row 136
$ws.Range("H136").Value = 22215.416
$ws.Range("I136").Value = 1870.3
$ws.Range("J136").Value = 36747.645
$ws.Range("K136").Value = 5610.9
$ws.Range("L136").Value = 110242.935
$ws.Range("M136").Value = -3060.9
$ws.Range("N136").Value = -115342.935

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1227.75
$ws.Range("J5").Value = 1407
$ws.Range("L5").Value = 4221
$ws.Range("N5").Value = -4445
# Row 92
$ws.Range("H92").Value = 750
$ws.Range("I92").Value = 366.66666
$ws.Range("K92").Value = 1099.99998
$ws.Range("M92").Value = 148.0000199999999
# Row 122
$ws.Range("H122").Value = 427
$ws.Range("I122").Value = 427
$ws.Range("K122").Value = 3843
$ws.Range("M122").Value = -1393
# Row 131
$ws.Range("H131").Value = 756.62
$ws.Range("J131").Value = 777.81055
$ws.Range("L131").Value = 2333.43165
$ws.Range("N131").Value = -12413.43165
# Row 135
$ws.Range("H135").Value = 1227.75
$ws.Range("J135").Value = 1407
$ws.Range("L135").Value = 12663
$ws.Range("N135").Value = -17733

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 2315.5
$ws.Range("I43").Value = 2315.5
$ws.Range("K43").Value = 2315.5
$ws.Range("M43").Value = -2164.5
# Row 46
$ws.Range("H46").Value = 22650
$ws.Range("J46").Value = 22650
$ws.Range("L46").Value = 22650
$ws.Range("N46").Value = -22962

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 1178571.4
$ws.Range("I2").Value = 1230769.2
$ws.Range("K2").Value = 1230769.2
$ws.Range("M2").Value = -1230657.2
# Row 132
$ws.Range("H132").Value = 2730.6155
$ws.Range("I132").Value = 1912.5
$ws.Range("J132").Value = 4039.6
$ws.Range("K132").Value = 5737.5
$ws.Range("L132").Value = 12118.8
$ws.Range("M132").Value = -3207.5
$ws.Range("N132").Value = -17178.8

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
